$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1 from "Point" to "Fullname"
$ws.Range("B1").Value = "Fullname"

# Move active selection to B2 (was B3)
$ws.Range("B2").Select()
